# "Sheet1" is the data-generator sheet: A2:C21 / A23:C42 are CONCATENATE()
# formulas seeded off I2 / I23. Bumping the seed numbers regenerates all of
# the test usernames/emails used elsewhere in the workbook.
$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("login")
$wsOrder = $wb.Worksheets.Item("order")
$wsData  = $wb.Worksheets.Item("Sheet1")

$wsData.Range("I2").Value = 21
$wsData.Range("I23").Value = 26
$excel.Calculate()

# "login" sheet columns G:I (rows 2-21) are a paste-as-values copy of
# Sheet1!A2:C21 - refresh them from the recalculated generator sheet.
$wsLogin.Range("G2:I21").Value = $wsData.Range("A2:C21").Value()

# "order" sheet columns R:T (rows 2-21) are a paste-as-values copy of
# Sheet1!A23:C42 - refresh them too.
$wsOrder.Range("R2:T21").Value = $wsData.Range("A23:C42").Value()

# Reflect the new selection on the "order" sheet that comes from reviewing
# the refreshed R:T block.
$wsOrder.Activate()
$wsOrder.Range("R2:T21").Select()
